$d = $word.ActiveDocument

# 1. Split " sensor,tmp sensor with below functionalities: " text.
$d.Content.Find.Execute("sensor,tmp sensor with below functionalities:", $true, $false, $false, $false, $false, $true, 1, $false, "sensor,tmp sensor with below functionalities:", 2)

Write-Output "done"
